$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to text format so numeric-looking strings are not
# auto-converted to numbers when we assign .Value below.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "34.410.12"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "1.803.49"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "227.80"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").Value = "0.581"
$ws.Range("E6").Value = "  +4.06%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "35.02"
$ws.Range("E8").Value = "  +6.35%  "
$ws.Range("D9").Value = "0.300"
$ws.Range("E9").Value = "  +0.69%  "
$ws.Range("D10").Value = "0.0692"
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").Value = "2.063.77"
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "11.17"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.801.30"
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("D15").Value = "0.641"
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("D16").Value = "34.384.75"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("D17").Value = "4.34"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("D18").Value = "68.98"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").Value = "245.16"
$ws.Range("E19").Value = "  -1.28%  "
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("E21").Value = "  +1.83%  "
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("E23").Value = "  -0.82%  "
$ws.Range("D24").Value = "170.39"
$ws.Range("E24").Value = "  +2.99%  "
$ws.Range("D25").Value = "2.11"
$ws.Range("E25").Value = "  +2.65%  "
$ws.Range("D26").Value = "7.57"
$ws.Range("E26").Value = "  +3.97%  "
$ws.Range("D27").Value = "16.71"
$ws.Range("E27").Value = "  +0.74%  "
$ws.Range("D28").Value = "0.119"
$ws.Range("E28").Value = "  +1.89%  "
$ws.Range("E29").Value = "  -0.58%  "
$ws.Range("D30").Value = "3.98"
$ws.Range("E30").Value = "  -4.89%  "
$ws.Range("D31").Value = "0.0528"
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("E32").Value = "  +0.70%  "
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "1.396.41"
$ws.Range("E35").Value = "  -1.92%  "
$ws.Range("D36").Value = "0.680"
$ws.Range("E36").Value = "  +0.72%  "
$ws.Range("D37").Value = "2.53"
$ws.Range("E37").Value = "  -2.71%  "
$ws.Range("E39").Value = "  -1.23%  "
$ws.Range("D40").Value = "82.97"
$ws.Range("E40").Value = "  -3.07%  "
$ws.Range("D41").Value = "2.84"
$ws.Range("E41").Value = "  +2.97%  "
$ws.Range("D42").Value = "0.946"
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("D44").Value = "13.56"
$ws.Range("E44").Value = "  -0.28%  "
$ws.Range("E45").Value = "  +2.68%  "
$ws.Range("D46").Value = "0.0510"
$ws.Range("E46").Value = "  -2.92%  "
$ws.Range("E47").Value = "  -1.39%  "
$ws.Range("D48").Value = "1.963.91"
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("D49").Value = "104.61"
$ws.Range("E49").Value = "  -1.60%  "
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("E51").Value = "  +0.97%  "

# Restore default (Normal) style so the cells keep the same
# style index as before (no explicit s="..." attribute).
$dataRange.Style = "Normal"
